# Weekly update: insert two new price rows (newest week) ahead of the
# existing history, pushing the previous rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 86 so the existing
# rows 86:90 shift down to 88:92 (carrying their formatting with them).
$ws.Rows("86:87").Insert()

# New row 86 - latest "Primera" quality entry ($/caja 18 kilos).
$ws.Range("A86").Value() = 9
$ws.Range("B86").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C86").Value() = "Metropolitana"
$ws.Range("D86").Value() = 44753
$ws.Range("E86").Value() = 13
$ws.Range("F86").Value() = 100114002
$ws.Range("G86").Value() = "Camote"
$ws.Range("H86").Value() = "Sin especificar"
$ws.Range("I86").Value() = "Primera"
$ws.Range("J86").Value() = 520
$ws.Range("K86").Value() = 12000
$ws.Range("L86").Value() = 13000
$ws.Range("M86").Value() = 12500
$ws.Range("N86").Value() = "`$/caja 18 kilos"
$ws.Range("O86").Value() = "Perú"
$ws.Range("P86").Value() = 694
$ws.Range("Q86").Value() = 18
$ws.Range("R86").Value() = "Hortaliza"

# New row 87 - latest "Primera" quality entry ($/malla 18 kilos).
$ws.Range("A87").Value() = 9
$ws.Range("B87").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value() = "Metropolitana"
$ws.Range("D87").Value() = 44753
$ws.Range("E87").Value() = 13
$ws.Range("F87").Value() = 100114002
$ws.Range("G87").Value() = "Camote"
$ws.Range("H87").Value() = "Sin especificar"
$ws.Range("I87").Value() = "Primera"
$ws.Range("J87").Value() = 1060
$ws.Range("K87").Value() = 9000
$ws.Range("L87").Value() = 10000
$ws.Range("M87").Value() = 9500
$ws.Range("N87").Value() = "`$/malla 18 kilos"
$ws.Range("O87").Value() = "Perú"
$ws.Range("P87").Value() = 528
$ws.Range("Q87").Value() = 18
$ws.Range("R87").Value() = "Hortaliza"

# Ensure the date cells keep the existing date/time number format.
$ws.Range("D86:D87").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
